$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = "T3uc"
$ws.Range("H3").Value = "T3uc"
$ws.Range("H4").Value = "T3uc"
$ws.Range("H5").Value = "T3uc"
$ws.Range("H6").Value = "ft2A"
$ws.Range("H7").Value = "01qd"
$ws.Range("H8").Value = "JgI1"
$ws.Range("H9").Value = "A6sB"
$ws.Range("H10").Value = "A6sB"
$ws.Range("H11").Value = "NYfe"
$ws.Range("H12").Value = "klkQ"
$ws.Range("H13").Value = "raQm"
$ws.Range("H14").Value = "gN8p"
$ws.Range("H15").Value = "ow6Q"
$ws.Range("H16").Value = "FPzC"
$ws.Range("H17").Value = "1uDu"
$ws.Range("H18").Value = "9B2h"
$ws.Range("H19").Value = "u4cY"
$ws.Range("H20").Value = "uqeM"
$ws.Range("H21").Value = "GoeX"
$ws.Range("H22").Value = "Khu3"
$ws.Range("H23").Value = "R7AU"
$ws.Range("H24").Value = "W6Zd"
$ws.Range("H25").Value = "gnYc"
$ws.Range("H26").Value = "jK4V"
$ws.Range("H27").Value = "jK4V"
$ws.Range("H28").Value = "caEi"
$ws.Range("H29").Value = "KQK8"
$ws.Range("H30").Value = "g4xo"
$ws.Range("H31").Value = "4HcB"
$ws.Range("H32").Value = "qqC2"
$ws.Range("H33").Value = "4F7R"
$ws.Range("H34").Value = "4F7R"
$ws.Range("H35").Value = "CN6X"
$ws.Range("H36").Value = "Tgl5"
$ws.Range("H37").Value = "yUWm"
$ws.Range("H38").Value = "fKPQ"
$ws.Range("H39").Value = "HCeS"
$ws.Range("H40").Value = "HCeS"
$ws.Range("H41").Value = "GuZL"
$ws.Range("H42").Value = "GuZL"
$ws.Range("H43").Value = "SCVF"
$ws.Range("H44").Value = "SCVF"
$ws.Range("H45").Value = "t5zw"
$ws.Range("H46").Value = "t5zw"
$ws.Range("H47").Value = "Sxxm"
$ws.Range("H48").Value = "Sxxm"
$ws.Range("H49").Value = "gNdd"
$ws.Range("H50").Value = "gNdd"
$ws.Range("H51").Value = "8val"
$ws.Range("H52").Value = "ww9m"
$ws.Range("H53").Value = "ww9m"
$ws.Range("H54").Value = "awPQ"
$ws.Range("H55").Value = "awPQ"
$ws.Range("H56").Value = "1IAR"
$ws.Range("H57").Value = "1IAR"
$ws.Range("H58").Value = "1IAR"
$ws.Range("H59").Value = "WGUu"
$ws.Range("H60").Value = "WGUu"
$ws.Range("H61").Value = "39dd"
$ws.Range("H62").Value = "OCyE"
$ws.Range("H63").Value = "S1Hc"
$ws.Range("H64").Value = "swXk"
$ws.Range("H65").Value = "JsLv"
$ws.Range("H66").Value = "jgJ7"
$ws.Range("H67").Value = "xLT2"
$ws.Range("H68").Value = "uJzC"
$ws.Range("H69").Value = "znju"
$ws.Range("H70").Value = "Q1Hl"
$ws.Range("H71").Value = "osaF"
$ws.Range("H72").Value = "IdhR"
$ws.Range("H73").Value = "gbvI"
$ws.Range("H74").Value = "LpXa"
$ws.Range("H75").Value = "iv40"
$ws.Range("H76").Value = "nTV5"
